# CRM.xlsx maintenance edit:
#  - insert a blank row above the header row
#  - fix header typos (Area -> Área, Fecha de Creaciòn -> Fecha de Creación,
#    Presupuesto Liquidacion -> Presupuesto Liquidación, Variacón +/- -> Variación +/-)
#  - apply a date number format to the "Fecha de Creación" / "Fecha Liquidación" header cells
#  - widen a number of data columns
#  - add a small accent/marker formatting in the new blank row (A1/B1)
#  - switch the page to portrait orientation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new blank row 1 - shifts the header row from row 1 to row 2
$ws.Rows.Item(1).Insert()

# --- fix header typos (this also renames the shared-string entries in place) ---
$ws.Cells.Item(2, 28).Value = "Variación +"
$ws.Cells.Item(2, 29).Value = "Variación -"
$ws.Cells.Item(2, 18).Value = "Presupuesto Liquidación"
$ws.Cells.Item(2, 3).Value  = "Área"
$ws.Cells.Item(2, 10).Value = "Fecha de Creación"

# --- date format on the "Fecha de Creación" (J) and "Fecha Liquidación" (N) headers ---
$dateFmt = "[$-10C6B]dd/mm/yyyy;@"
$ws.Cells.Item(2, 10).NumberFormat = $dateFmt
$ws.Cells.Item(2, 14).NumberFormat = $dateFmt

# --- widen columns ---
$offset = 0.8333333333333334
$ws.Columns.Item(15).ColumnWidth = 15.26171875 - $offset
$ws.Columns.Item(16).ColumnWidth = 13.15625 - $offset
$ws.Columns.Item(17).ColumnWidth = 14.3125 - $offset
$ws.Columns.Item(18).ColumnWidth = 16.68359375 - $offset
$ws.Columns.Item(19).ColumnWidth = 15.41796875 - $offset
$ws.Columns.Item(20).ColumnWidth = 14.9453125 - $offset
$ws.Columns.Item(24).ColumnWidth = 18.68359375 - $offset
$ws.Columns.Item(28).ColumnWidth = 14.3671875 - $offset
$ws.Columns.Item(29).ColumnWidth = 16.89453125 - $offset
$ws.Columns.Item(30).ColumnWidth = 14.83984375 - $offset

# --- small marker cells in the new blank row ---
$ws.Cells.Item(1, 2).NumberFormat = $dateFmt
$ws.Cells.Item(1, 1).NumberFormat = '_ "S/."\ * #,##0.00_ ;_ "S/."\ * \-#,##0.00_ ;_ "S/."\ * "-"??_ ;_ @_ '
$ws.Cells.Item(1, 1).Font.Name = "Calibri"
$ws.Cells.Item(1, 1).Font.Size = 11
$ws.Cells.Item(1, 1).Interior.ThemeColor = 0
$ws.Cells.Item(1, 1).Interior.Pattern = 1
$ws.Cells.Item(1, 1).Borders.LineStyle = 4
$ws.Cells.Item(1, 1).HorizontalAlignment = -4108
$ws.Cells.Item(1, 1).VerticalAlignment = -4108

# --- page orientation ---
$ws.PageSetup.Orientation = 1

Write-Output "done"
